# "Generate Report for Handoff"
#
# Updates the localization-status report to reflect that the content is
# now Ready for handoff (was "In Translation"), refreshes the
# handoff/generation timestamps, and widens the Status columns so the
# new, longer status text fits (mirrors Excel's column auto-fit).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value2 = "Ready for handoff"
$wsOverview.Range("F2").Value2 = "Ready for handoff"
$wsZhCn.Range("C2").Value2     = "Ready for handoff"
$wsDeDe.Range("C2").Value2     = "Ready for handoff"

# --- Refresh timestamps ----------------------------------------------------
# Overview: "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value2 = "2016-08-24 20:40:31"

# zh-cn: "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value2 = "2016-08-24 20:40:27"

# de-de: "Latest Handoff Datetime"
$wsDeDe.Range("H2").Value2 = "2016-08-24 20:40:31"

# --- Widen the Status columns to fit "Ready for handoff" -------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33   # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33   # C: Status
